$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.295.70"
$ws.Range("E2").Value = "  -0.07%  "
# Row 3
$ws.Range("D3").Value = "1.927.80"
$ws.Range("E3").Value = "  -0.21%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7476"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.86%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.03%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3147"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.58%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06950"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.42%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07999"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.15%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7691"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.74%  "
# Row 13
$ws.Range("D13").Value = "1.925.97"
$ws.Range("E13").Value = "  -0.31%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.300"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.39%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.90%  "
# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "30.297.68"
$ws.Range("E16").Value = "  -0.08%  "
# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.80%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "250.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007872"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.97%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.713"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "
# Row 21
$ws.Range("D21").Value = "2.186.29"
$ws.Range("E21").Value = "  +0.35%  "
# Row 22
$ws.Range("E22").Value = "  +0.13%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.615"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.38%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.398"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.11%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1321"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.01%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.183"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.90%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.368"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.59%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.511"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.353"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.71%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.082"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.84%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.271"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7410"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.780"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.59%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01941"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.12%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "
# Row 40
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.93%  "
# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.349"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.24%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4411"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.58%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.946"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.05%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8312"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.52%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.666"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.413"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.06%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "969.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.35%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06035"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
